# Applies the data refresh described in the commit:
# "Updated cryptos list on Sat Nov 30 12:00:13 UTC 2024 with GitHub Actions"
# Each row's Price (D) / Volume(1h) (E) values are refreshed; a handful of rows
# were additionally re-ranked, swapping the Coin name (B) and Link (C) between
# adjacent rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '96.320.19'
$ws.Range("D3").Value = '3.658.32'
$ws.Range("E3").Value = '  +1.72%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '241.59'
$ws.Range("E5").Value = '  -1.48%  '
$ws.Range("E6").Value = '  +14.87%  '
$ws.Range("D7").Value = '659.52'
$ws.Range("E7").Value = '  +0.47%  '
$ws.Range("D8").Value = '0.422'
$ws.Range("E8").Value = '  +2.13%  '
$ws.Range("E9").Value = '  +2.82%  '
$ws.Range("D11").Value = '3.657.06'
$ws.Range("E11").Value = '  +1.76%  '
$ws.Range("D12").Value = '44.71'
$ws.Range("E12").Value = '  +2.09%  '
$ws.Range("E13").Value = '  +0.17%  '
$ws.Range("D14").Value = '6.65'
$ws.Range("E14").Value = '  +3.04%  '
$ws.Range("D15").Value = '4.341.47'
$ws.Range("D16").Value = '0.0000272'
$ws.Range("E16").Value = '  +5.30%  '
$ws.Range("D17").Value = '96.088.44'
$ws.Range("E17").Value = '  -1.06%  '
$ws.Range("D18").Value = '8.90'
$ws.Range("E18").Value = '  +14.92%  '
$ws.Range("D19").Value = '3.656.30'
$ws.Range("E19").Value = '  +2.64%  '
$ws.Range("D20").Value = '12.70'
$ws.Range("E20").Value = '  -0.34%  '
$ws.Range("D21").Value = '18.23'
$ws.Range("E21").Value = '  +0.64%  '
$ws.Range("D22").Value = '0.532'
$ws.Range("D23").Value = '519.04'
$ws.Range("E23").Value = '  +1.54%  '
$ws.Range("D24").Value = '3.42'
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").Value = '0.0000203'
$ws.Range("E25").Value = '  +0.82%  '
$ws.Range("D26").Value = '6.86'
$ws.Range("E26").Value = '  -1.00%  '
$ws.Range("D27").Value = '101.96'
$ws.Range("E27").Value = '  +4.76%  '
$ws.Range("D28").Value = '12.93'
$ws.Range("E28").Value = '  -1.98%  '
$ws.Range("D29").Value = '0.168'
$ws.Range("E29").Value = '  +11.10%  '
$ws.Range("D30").Value = '12.35'
$ws.Range("E30").Value = '  +6.47%  '
$ws.Range("D31").Value = '3.03'
$ws.Range("E31").Value = '  -0.89%  '
$ws.Range("E32").Value = '  +0.12%  '
$ws.Range("E33").Value = '  -0.14%  '
$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").Value = '1.01'
$ws.Range("E34").Value = '  +1.55%  '
$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").Value = '1.83'
$ws.Range("E35").Value = '  +11.57%  '
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").Value = '32.81'
$ws.Range("E36").Value = '  +3.68%  '
$ws.Range("E37").Value = '  +2.32%  '
$ws.Range("D38").Value = '624.35'
$ws.Range("E38").Value = '  +0.08%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").Value = '8.68'
$ws.Range("E39").Value = '  -1.67%  '
$ws.Range("B40").Value = 'EnergySwap'
$ws.Range("C40").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D40").Value = '44.39'
$ws.Range("E40").Value = '  +33.23%  '
$ws.Range("E41").Value = '  +4.40%  '
$ws.Range("D42").Value = '0.954'
$ws.Range("E42").Value = '  +4.06%  '
$ws.Range("E43").Value = '  +4.63%  '
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("E45").Value = '  +7.43%  '
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = '0.452'
$ws.Range("E46").Value = '  +25.80%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").Value = '0.0457'
$ws.Range("E47").Value = '  +5.18%  '
$ws.Range("B48").Value = 'WhiteBITCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D48").Value = '23.61'
$ws.Range("E48").Value = '  -0.15%  '
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").Value = '2.27'
$ws.Range("E49").Value = '  -1.63%  '
$ws.Range("D50").Value = '8.54'
$ws.Range("D51").Value = '3.56'
$ws.Range("E51").Value = '  +1.03%  '
